# Add MySQL connect info: insert 4 new columns (SqlPort, SqlName, SqlUser, SqlPwd
# area) right before the old "F" column, then populate the new header/data cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 4 columns before the old column F (shifts old F:J -> J:N, along with
#    their styles / data validation / column widths).
$ws.Range("F1:I1").EntireColumn.Insert()

# 2) Re-purpose the old "Pwd" header (E1) as "SqlIP", and fill in the new headers.
$ws.Range("E1").Value = "SqlIP"
$ws.Range("F1").Value = "SqlPort"
$ws.Range("G1").Value = "SqlName"
$ws.Range("H1").Value = "SqlUser"
$ws.Range("I1").Value = "SqlPwd"

# 3) Fill in the MySQL connection info row.
$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("F2").Value = 3306
$ws.Range("G2").Value = "app_test"
$ws.Range("H2").Value = "root"
$ws.Range("I2").Value = 123456

# 4) Column widths (best fit for the new layout).
$ws.Columns.Item(1).ColumnWidth = 12.035714285714286
$ws.Columns.Item(4).ColumnWidth = 4.785714285714286
$ws.Columns.Item(5).ColumnWidth = 13.160714285714286
$ws.Columns.Item(6).ColumnWidth = 7.785714285714286
$ws.Columns.Item(7).ColumnWidth = 7.785714285714286
$ws.Columns.Item(8).ColumnWidth = 7.785714285714286
$ws.Columns.Item(9).ColumnWidth = 6.785714285714286

# 5) Selection moves to G9.
$null = $ws.Range("G9").Select()
